$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 35, pushing existing rows 35-38 (STATE, ERROR, QS_TX, QS_RX) down to 36-39
$ws.Rows.Item(35).Insert()

# Populate the new row 35 with the WATER_DEPTH datapoint
$ws.Range("A35").Value = "datapoints"
$ws.Range("B35").Value = "WATER_DEPTH"
$ws.Range("C35").Value = "m"
$ws.Range("D35").Value = "Water depth"

$wb.Save()
